$wb = $excel.ActiveWorkbook

# ---- Query1 (sheet1) data update ----
$ws1 = $wb.Worksheets.Item("Query1")
$data1 = @(
    @(0.725,0.81,0.775,0.775,0.815,0.785,0.79,0.795,0.755,0.8),
    @(0.8365,0.739,0.7082,0.7854,0.7493,0.8192,0.732,0.7808,0.7734,0.7581),
    @(0.8,0.8,0.775,0.775,0.7775,0.7675,0.755,0.735,0.735,0.7325),
    @(0.7742,0.7561,0.7442,0.8032,0.7623,0.8058,0.7801,0.7006,0.7471,0.8218),
    @(0.7567,0.78,0.7833,0.7517,0.745,0.7633,0.7767,0.7883,0.78,0.7383),
    @(0.8104,0.7701,0.814,0.7755,0.8055,0.8173,0.7562,0.7384,0.7647,0.7307),
    @(0.7913,0.7475,0.7575,0.7863,0.7875,0.7888,0.7513,0.77,0.7775,0.7713),
    @(0.7795,0.7835,0.7926,0.7717,0.7755,0.8032,0.7798,0.8051,0.7604,0.7568),
    @(0.742,0.743,0.789,0.769,0.782,0.758,0.782,0.793,0.772,0.78),
    @(0.7523,0.7682,0.7732,0.793,0.7813,0.7474,0.757,0.769,0.7806,0.7671)
)

for ($r = 0; $r -lt $data1.Length; $r++) {
    $rowVals = $data1[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws1.Cells.Item($r + 2, $c + 3).Value = $rowVals[$c]
    }
}

# ---- Query2 (sheet2) data update ----
$ws2 = $wb.Worksheets.Item("Query2")
$data2 = @(
    @(0.58,0.465,0.505,0.445,0.505,0.465,0.575,0.51,0.455,0.495),
    @(0.52,0.52,0.515,0.5,0.55,0.475,0.53,0.455,0.545,0.47),
    @(0.5075,0.4825,0.4825,0.475,0.4925,0.495,0.5175,0.485,0.475,0.54),
    @(0.4475,0.515,0.475,0.46,0.49,0.4925,0.485,0.4825,0.5175,0.4925),
    @(0.52,0.495,0.4783,0.505,0.497,0.5017,0.515,0.515,0.475,0.505),
    @(0.4817,0.5233,0.5033,0.515,0.497,0.4783,0.5217,0.4983,0.4917,0.5017),
    @(0.5313,0.4775,0.5125,0.5038,0.5413,0.5163,0.4863,0.5263,0.5113,0.505),
    @(0.51,0.5088,0.4787,0.5038,0.4825,0.495,0.5125,0.4787,0.51,0.4837),
    @(0.519,0.502,0.498,0.481,0.473,0.483,0.47,0.506,0.524,0.497),
    @(0.488,0.499,0.484,0.533,0.506,0.501,0.513,0.495,0.52,0.496)
)

for ($r = 0; $r -lt $data2.Length; $r++) {
    $rowVals = $data2[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws2.Cells.Item($r + 2, $c + 3).Value = $rowVals[$c]
    }
}

# ---- Selection state: both sheets now have A12 selected ----
$ws1.Range("A12").Select()
$ws2.Range("A12").Select()

# Query2 remains the active/selected tab (tabSelected="true"), matching original + target state
$ws2.Activate()
